$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "mahara" provider rows (139-149) now have BU values fetched into column D.
$buValues = @{
    139 = 4686
    140 = 4640
    141 = 4686
    142 = 562
    143 = 4681
    144 = 4681
    145 = 4640
    146 = 4642
    147 = 311
    148 = 4541
    149 = 4527
}

foreach ($row in $buValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $buValues[$row]
}

# The BU column (D) now carries numeric data for the whole table, so apply the
# same number-formatted style used by column C (Amt.) across all data rows.
$ws.Range("D2:D230").NumberFormat = "#,##0;[Red]-#,##0"
$ws.Range("D2:D230").VerticalAlignment = -4108
